# "Single Author Customized query" - append the new lookup result as row 9.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "Robert Schifreen"
$ws.Range("B9").Value = "Journals"

# Column A needs to be a bit wider to fit the new (longer) author name.
$ws.Columns("A").ColumnWidth = 30.3

# Leave the selection where the author left it when they saved the file.
$ws.Range("C10").Select()
